$wb = $excel.ActiveWorkbook

# Sheet "OFF" - row 2 (Home) updates
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 167
$wsOff.Range("C2").Value = 106
$wsOff.Range("D2").Value = 42
$wsOff.Range("E2").Value = 23

# Sheet "DEF" - row 2 (Home) updates
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 154
$wsDef.Range("C2").Value = 110
$wsDef.Range("D2").Value = 39
$wsDef.Range("G2").Value = 2
